$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")
$ws.Range("T47").Formula = "=S47*0.7"
